$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (from original row 35)
$ws.Range("D2").Value = 44874
$ws.Range("J2").Value = 7900

# Row 3 (from original row 21)
$ws.Range("D3").Value = 44162
$ws.Range("J3").Value = 7000
$ws.Range("O3").Value = "Provincia de Chacabuco"

# Row 4 (from original row 38)
$ws.Range("D4").Value = 44232
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 16000
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 3000
$ws.Range("O4").Value = "Provincia de Chacabuco"
$ws.Range("P4").Value = 30

# Row 5 (from original row 37)
$ws.Range("D5").Value = 44166
$ws.Range("J5").Value = 7000

# Row 6 (from original row 41)
$ws.Range("D6").Value = 44181
$ws.Range("J6").Value = 12000

# Row 8 (from original row 20)
$ws.Range("D8").Value = 44215
$ws.Range("J8").Value = 16000

# Row 9 (from original row 26)
$ws.Range("D9").Value = 44161
$ws.Range("J9").Value = 7000

# Row 10 (from original row 28)
$ws.Range("D10").Value = 44845

# Row 11 (from original row 39)
$ws.Range("D11").Value = 44901
$ws.Range("J11").Value = 7000

# Row 12 (from original row 36)
$ws.Range("D12").Value = 44873

# Row 13 (from original row 3)
$ws.Range("D13").Value = 44245
$ws.Range("J13").Value = 9000
$ws.Range("O13").Value = "Región Metropolitana"

# Row 14 (from original row 4)
$ws.Range("D14").Value = 44245
$ws.Range("I14").Value = "Segunda"
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 2500
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2500
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 25

# Row 15 (from original row 40)
$ws.Range("D15").Value = 44188
$ws.Range("J15").Value = 12000

# Row 16 (from original row 12)
$ws.Range("D16").Value = 44882

# Row 17 (from original row 31)
$ws.Range("D17").Value = 44859
$ws.Range("J17").Value = 7900

# Row 18 (from original row 32)
$ws.Range("D18").Value = 44602
$ws.Range("J18").Value = 12000

# Row 19 (from original row 33)
$ws.Range("D19").Value = 44602
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value = 6000
$ws.Range("K19").Value = 2500
$ws.Range("L19").Value = 2500
$ws.Range("M19").Value = 2500
$ws.Range("P19").Value = 25

# Row 20 (from original row 29)
$ws.Range("D20").Value = 44855
$ws.Range("J20").Value = 7900

# Row 21 (from original row 23)
$ws.Range("D21").Value = 44160

# Row 22 (from original row 15)
$ws.Range("D22").Value = 44167
$ws.Range("J22").Value = 7000

# Row 23 (from original row 22)
$ws.Range("D23").Value = 44186
$ws.Range("J23").Value = 10000

# Row 24 (from original row 25)
$ws.Range("D24").Value = 44600
$ws.Range("J24").Value = 1300
$ws.Range("K24").Value = 3500
$ws.Range("L24").Value = 4000
$ws.Range("M24").Value = 3808
$ws.Range("O24").Value = "Región Metropolitana"
$ws.Range("P24").Value = 38

# Row 25 (from original row 30)
$ws.Range("D25").Value = 44847
$ws.Range("J25").Value = 7900
$ws.Range("K25").Value = 3000
$ws.Range("L25").Value = 3000
$ws.Range("M25").Value = 3000
$ws.Range("O25").Value = "Provincia de Chacabuco"
$ws.Range("P25").Value = 30

# Row 26 (from original row 16)
$ws.Range("D26").Value = 44876
$ws.Range("J26").Value = 7900

# Row 27 (from original row 34)
$ws.Range("D27").Value = 44210
$ws.Range("J27").Value = 8800
$ws.Range("K27").Value = 2500
$ws.Range("M27").Value = 2750
$ws.Range("P27").Value = 28

# Row 28 (from original row 18)
$ws.Range("D28").Value = 44159
$ws.Range("J28").Value = 7000

# Row 29 (from original row 43)
$ws.Range("D29").Value = 44880

# Row 30 (from original row 17)
$ws.Range("D30").Value = 44229
$ws.Range("J30").Value = 16000

# Row 31 (from original row 10)
$ws.Range("D31").Value = 44860

# Row 32 (from original row 2)
$ws.Range("D32").Value = 44231

# Row 33 (from original row 13)
$ws.Range("D33").Value = 44187
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 12000
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 3000
$ws.Range("M33").Value = 3000
$ws.Range("P33").Value = 30

# Row 34 (from original row 8)
$ws.Range("D34").Value = 44902
$ws.Range("J34").Value = 7000
$ws.Range("K34").Value = 3000
$ws.Range("M34").Value = 3000
$ws.Range("P34").Value = 30

# Row 35 (from original row 19)
$ws.Range("D35").Value = 44883
$ws.Range("J35").Value = 9700

# Row 36 (from original row 11)
$ws.Range("D36").Value = 44189
$ws.Range("J36").Value = 16000

# Row 37 (from original row 6)
$ws.Range("D37").Value = 44230
$ws.Range("J37").Value = 16000

# Row 38 (from original row 5)
$ws.Range("D38").Value = 44875
$ws.Range("J38").Value = 7900

# Row 39 (from original row 42)
$ws.Range("D39").Value = 44209
$ws.Range("K39").Value = 2500
$ws.Range("M39").Value = 2750
$ws.Range("P39").Value = 28

# Row 40 (from original row 9)
$ws.Range("D40").Value = 44846
$ws.Range("J40").Value = 7900

# Row 41 (from original row 24)
$ws.Range("D41").Value = 44204
$ws.Range("J41").Value = 7000

# Row 42 (from original row 27)
$ws.Range("D42").Value = 44214
$ws.Range("K42").Value = 3000
$ws.Range("M42").Value = 3000
$ws.Range("P42").Value = 30

# Row 43 (from original row 14)
$ws.Range("D43").Value = 44881
